$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 328.4
$ws.Cells.Item(28, 9).Value = 178.41667
$ws.Cells.Item(28, 11).Value = 178.41667
$ws.Cells.Item(28, 13).Value = 306.58333
$ws.Cells.Item(41, 8).Value = 375
$ws.Cells.Item(41, 9).Value = 433.33334
$ws.Cells.Item(41, 10).Value = 287.5
$ws.Cells.Item(41, 11).Value = 433.33334
$ws.Cells.Item(41, 12).Value = 287.5
$ws.Cells.Item(41, 13).Value = 6.666659999999979
$ws.Cells.Item(41, 14).Value = -1167.5
$ws.Cells.Item(62, 8).Value = 3952.3225
$ws.Cells.Item(62, 9).Value = 3949.84
$ws.Cells.Item(62, 10).Value = 3962.6667
$ws.Cells.Item(62, 11).Value = 3949.84
$ws.Cells.Item(62, 12).Value = 3962.6667
$ws.Cells.Item(62, 13).Value = -3325.84
$ws.Cells.Item(62, 14).Value = -5210.6667
$ws.Cells.Item(65, 8).Value = 3952.3225
$ws.Cells.Item(65, 9).Value = 3949.84
$ws.Cells.Item(65, 10).Value = 3962.6667
$ws.Cells.Item(65, 11).Value = 19749.2
$ws.Cells.Item(65, 12).Value = 19813.3335
$ws.Cells.Item(65, 13).Value = -16629.2
$ws.Cells.Item(65, 14).Value = -26053.3335
$ws.Cells.Item(82, 8).Value = 1399.2
$ws.Cells.Item(82, 9).Value = 1399.2
$ws.Cells.Item(82, 11).Value = 4197.6
$ws.Cells.Item(82, 13).Value = -3791.6
$ws.Cells.Item(85, 8).Value = 1399.2
$ws.Cells.Item(85, 9).Value = 1399.2
$ws.Cells.Item(85, 11).Value = 4197.6
$ws.Cells.Item(85, 13).Value = -2793.6
$ws.Cells.Item(96, 8).Value = 467.83334
$ws.Cells.Item(96, 9).Value = 467.83334
$ws.Cells.Item(96, 11).Value = 1403.50002
$ws.Cells.Item(96, 13).Value = -30.50001999999995
$ws.Cells.Item(106, 8).Value = 3227.4
$ws.Cells.Item(106, 9).Value = 2468.4285
$ws.Cells.Item(106, 10).Value = 4998.3335
$ws.Cells.Item(106, 11).Value = 2468.4285
$ws.Cells.Item(106, 12).Value = 4998.3335
$ws.Cells.Item(106, 13).Value = -1837.4285
$ws.Cells.Item(106, 14).Value = -6260.3335
$ws.Cells.Item(132, 8).Value = 5028.2856
$ws.Cells.Item(132, 9).Value = 5650.826
$ws.Cells.Item(132, 11).Value = 16952.478
$ws.Cells.Item(132, 13).Value = -14422.478
$ws.Cells.Item(137, 8).Value = 2865
$ws.Cells.Item(137, 9).Value = 2865
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 8595
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -6045
$ws.Cells.Item(137, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 300
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 300
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -184
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 2410.6365
$ws.Cells.Item(32, 9).Value = 2173.4688
$ws.Cells.Item(32, 11).Value = 2173.4688
$ws.Cells.Item(32, 13).Value = -1886.4688
$ws.Cells.Item(61, 8).Value = 3793643
$ws.Cells.Item(61, 9).Value = 5380835.5
$ws.Cells.Item(61, 11).Value = 5380835.5
$ws.Cells.Item(61, 13).Value = -5380623.5
$ws.Cells.Item(76, 8).Value = 20000
$ws.Cells.Item(76, 10).Value = 20000
$ws.Cells.Item(76, 12).Value = 20000
$ws.Cells.Item(76, 14).Value = -20676
$ws.Cells.Item(79, 8).Value = 20000
$ws.Cells.Item(79, 10).Value = 20000
$ws.Cells.Item(79, 12).Value = 20000
$ws.Cells.Item(79, 14).Value = -22340
$ws.Cells.Item(132, 8).Value = 5550.5947
$ws.Cells.Item(132, 9).Value = 4569.6665
$ws.Cells.Item(132, 11).Value = 13708.9995
$ws.Cells.Item(132, 13).Value = -11178.9995
$ws.Cells.Item(136, 8).Value = 3793643
$ws.Cells.Item(136, 9).Value = 5380835.5
$ws.Cells.Item(136, 11).Value = 16142506.5
$ws.Cells.Item(136, 13).Value = -16139956.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1940.6923
$ws.Cells.Item(20, 9).Value = 2144.5
$ws.Cells.Item(20, 10).Value = 1261.3334
$ws.Cells.Item(20, 11).Value = 2144.5
$ws.Cells.Item(20, 12).Value = 1261.3334
$ws.Cells.Item(20, 13).Value = -1897.5
$ws.Cells.Item(20, 14).Value = -1755.3334
$ws.Cells.Item(96, 8).Value = 10000
$ws.Cells.Item(96, 9).Value = 10000
$ws.Cells.Item(96, 11).Value = 10000
$ws.Cells.Item(96, 13).Value = -7254
$ws.Cells.Item(105, 8).Value = 7788.4443
$ws.Cells.Item(105, 9).Value = 7259.8
$ws.Cells.Item(105, 10).Value = 8449.25
$ws.Cells.Item(105, 11).Value = 7259.8
$ws.Cells.Item(105, 12).Value = 8449.25
$ws.Cells.Item(105, 13).Value = -5512.8
$ws.Cells.Item(105, 14).Value = -11943.25
$ws.Cells.Item(131, 8).Value = 99499.5
$ws.Cells.Item(131, 10).Value = 99499.5
$ws.Cells.Item(131, 12).Value = 99499.5
$ws.Cells.Item(131, 14).Value = -109579.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1853.75
$ws.Cells.Item(134, 8).Value = 3660.4075
$ws.Cells.Item(134, 9).Value = 1763.5238
$ws.Cells.Item(134, 10).Value = 10299.5
$ws.Cells.Item(134, 11).Value = 5290.5714
$ws.Cells.Item(134, 12).Value = 30898.5
$ws.Cells.Item(134, 13).Value = -2755.5714
$ws.Cells.Item(134, 14).Value = -35968.5
$ws.Cells.Item(137, 8).Value = 83750.664
$ws.Cells.Item(137, 10).Value = 83750.664
$ws.Cells.Item(137, 12).Value = 83750.664
$ws.Cells.Item(137, 14).Value = -93950.664
$ws.Cells.Item(141, 8).Value = 37299.5
$ws.Cells.Item(141, 10).Value = 36900
$ws.Cells.Item(141, 12).Value = 36900
$ws.Cells.Item(141, 14).Value = -47260

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(28, 8).Value = 2510
$ws.Cells.Item(28, 9).Value = 2765
$ws.Cells.Item(28, 10).Value = 2000
$ws.Cells.Item(28, 11).Value = 8295
$ws.Cells.Item(28, 12).Value = 6000
$ws.Cells.Item(28, 13).Value = -8063
$ws.Cells.Item(28, 14).Value = -6464
$ws.Cells.Item(37, 8).Value = 55920.617
$ws.Cells.Item(37, 10).Value = 55920.617
$ws.Cells.Item(37, 12).Value = 167761.851
$ws.Cells.Item(37, 14).Value = -167985.851
$ws.Cells.Item(131, 8).Value = 23078472
$ws.Cells.Item(131, 10).Value = 8336184.5
$ws.Cells.Item(131, 12).Value = 25008553.5
$ws.Cells.Item(131, 14).Value = -25018633.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 12777.5
$ws.Cells.Item(21, 10).Value = 15555
$ws.Cells.Item(21, 12).Value = 15555
$ws.Cells.Item(21, 14).Value = -15901
$ws.Cells.Item(30, 8).Value = 12777.5
$ws.Cells.Item(30, 10).Value = 15555
$ws.Cells.Item(30, 12).Value = 15555
$ws.Cells.Item(30, 14).Value = -15765
$ws.Cells.Item(47, 8).Value = 24992.334
$ws.Cells.Item(47, 9).Value = 24994.5
$ws.Cells.Item(47, 10).Value = 24988
$ws.Cells.Item(47, 11).Value = 24994.5
$ws.Cells.Item(47, 12).Value = 24988
$ws.Cells.Item(47, 13).Value = -24426.5
$ws.Cells.Item(47, 14).Value = -26124
$ws.Cells.Item(102, 8).Value = 2503
$ws.Cells.Item(102, 9).Value = 2297.5833
$ws.Cells.Item(102, 11).Value = 2297.5833
$ws.Cells.Item(102, 13).Value = -675.5832999999998
$ws.Cells.Item(107, 8).Value = 1073884.8
$ws.Cells.Item(107, 9).Value = 1825274.5
$ws.Cells.Item(107, 10).Value = 470.7143
$ws.Cells.Item(107, 11).Value = 1825274.5
$ws.Cells.Item(107, 12).Value = 470.7143
$ws.Cells.Item(107, 13).Value = -1823354.5
$ws.Cells.Item(107, 14).Value = -4310.7143
$ws.Cells.Item(122, 8).Value = 1826.7693
$ws.Cells.Item(122, 9).Value = 1870.6666
$ws.Cells.Item(122, 10).Value = 1300
$ws.Cells.Item(122, 11).Value = 5611.9998
$ws.Cells.Item(122, 12).Value = 3900
$ws.Cells.Item(122, 13).Value = -3161.9998
$ws.Cells.Item(122, 14).Value = -8800
$ws.Cells.Item(132, 8).Value = 5050.241
$ws.Cells.Item(132, 9).Value = 4658.44
$ws.Cells.Item(132, 10).Value = 7499
$ws.Cells.Item(132, 11).Value = 13975.32
$ws.Cells.Item(132, 12).Value = 22497
$ws.Cells.Item(132, 13).Value = -11445.32
$ws.Cells.Item(132, 14).Value = -27557

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 783.3333
$ws.Cells.Item(13, 9).Value = 566.6667
$ws.Cells.Item(13, 10).Value = 1000
$ws.Cells.Item(13, 11).Value = 566.6667
$ws.Cells.Item(13, 12).Value = 1000
$ws.Cells.Item(13, 13).Value = -426.6667
$ws.Cells.Item(13, 14).Value = -1280
$ws.Cells.Item(46, 8).Value = 23828080
$ws.Cells.Item(46, 10).Value = 35449.668
$ws.Cells.Item(46, 12).Value = 35449.668
$ws.Cells.Item(46, 14).Value = -35825.668
$ws.Cells.Item(93, 8).Value = 13409.526
$ws.Cells.Item(93, 9).Value = 821.9231
$ws.Cells.Item(93, 10).Value = 40682.668
$ws.Cells.Item(93, 11).Value = 821.9231
$ws.Cells.Item(93, 12).Value = 40682.668
$ws.Cells.Item(93, 13).Value = 426.0769
$ws.Cells.Item(93, 14).Value = -43178.668
$ws.Cells.Item(122, 8).Value = 4259.364
$ws.Cells.Item(122, 9).Value = 2539.2222
$ws.Cells.Item(122, 11).Value = 7617.6666
$ws.Cells.Item(122, 13).Value = -5167.6666
$ws.Cells.Item(132, 8).Value = 9051.210999999999
$ws.Cells.Item(132, 9).Value = 9835.286
$ws.Cells.Item(132, 11).Value = 29505.858
$ws.Cells.Item(132, 13).Value = -26975.858
$ws.Cells.Item(136, 8).Value = 4281.9697
$ws.Cells.Item(136, 9).Value = 3752.037
$ws.Cells.Item(136, 10).Value = 6666.6665
$ws.Cells.Item(136, 11).Value = 11256.111
$ws.Cells.Item(136, 12).Value = 19999.9995
$ws.Cells.Item(136, 13).Value = -8706.110999999999
$ws.Cells.Item(136, 14).Value = -25099.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4776.619
$ws.Cells.Item(132, 9).Value = 4200.7188
$ws.Cells.Item(132, 11).Value = 12602.1564
$ws.Cells.Item(132, 13).Value = -10072.1564
